$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, pushing existing rows 83:88 down to 84:89
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new data record
$ws.Cells.Item(83, 1).Value2  = 1
$ws.Cells.Item(83, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(83, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(83, 4).Value2  = 44748
$ws.Cells.Item(83, 5).Value2  = 15
$ws.Cells.Item(83, 6).Value2  = "Fruta"
$ws.Cells.Item(83, 7).Value2  = 100102
$ws.Cells.Item(83, 8).Value2  = "Cítricos"
$ws.Cells.Item(83, 9).Value2  = 100102005
$ws.Cells.Item(83, 10).Value2 = "Naranja"
$ws.Cells.Item(83, 11).Value2 = "Fukumoto"
$ws.Cells.Item(83, 12).Value2 = "Tercera"
$ws.Cells.Item(83, 13).Value2 = 300
$ws.Cells.Item(83, 14).Value2 = 500
$ws.Cells.Item(83, 15).Value2 = 600
$ws.Cells.Item(83, 16).Value2 = 550
$ws.Cells.Item(83, 17).Value2 = "`$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(83, 18).Value2 = "Región de Coquimbo"
$ws.Cells.Item(83, 19).Value2 = 550
$ws.Cells.Item(83, 20).Value2 = 1
